# Generate Report for handback
# - Marks the zh-cn and de-de handback rows as handed back (status, datetime)
# - Populates the "Latest Target File" / "Latest Handback File" columns (E/F)
#   with hyperlinks that mirror the existing source/handoff links
# - Sets the Handoff Reason to "Include" for the two localized files

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 - 5f0316cf-cd57-45e1-8eb3-c10d7ebc191b
$zh.Range("B2").Value = $newStatus
$zh.Range("E2").Value = "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md"
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0fc59f4a4e849b032212ae52abe8c2c84f5b9074/e2e/5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md", "", "", "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md")

$zh.Range("F2").Value = "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ad3789b038b0ae674e911ac73ff9c990ae4a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.zh-cn.xlf", "", "", "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.zh-cn.xlf")

$zh.Range("G2").Value = "2016-02-16 14:55:16"
$zh.Range("H2").Value = "Include"

# Row 3 - b06fe6dd-862b-45e4-b6c2-799b23eb6c96
$zh.Range("B3").Value = $newStatus
$zh.Range("E3").Value = "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md"
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0fc59f4a4e849b032212ae52abe8c2c84f5b9074/e2e/b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md", "", "", "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md")

$zh.Range("F3").Value = "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ad3789b038b0ae674e911ac73ff9c990ae4a9b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.zh-cn.xlf", "", "", "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.zh-cn.xlf")

$zh.Range("G3").Value = "2016-02-16 14:55:16"
$zh.Range("H3").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 - 5f0316cf-cd57-45e1-8eb3-c10d7ebc191b
$de.Range("B2").Value = $newStatus
$de.Range("E2").Value = "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md"
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/0fc59f4a4e849b032212ae52abe8c2c84f5b9074/e2e/5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md", "", "", "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.md")

$de.Range("F2").Value = "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.de-de.xlf"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5834b6e726c1c7b56b89edbab07a04b41dd1997a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.de-de.xlf", "", "", "5f0316cf-cd57-45e1-8eb3-c10d7ebc191b.1e3e331be7193c62f3327a90d13167a9e68b0e1e.de-de.xlf")

$de.Range("G2").Value = "2016-02-16 14:55:46"
$de.Range("H2").Value = "Include"

# Row 3 - b06fe6dd-862b-45e4-b6c2-799b23eb6c96
$de.Range("B3").Value = $newStatus
$de.Range("E3").Value = "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md"
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/0fc59f4a4e849b032212ae52abe8c2c84f5b9074/e2e/b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md", "", "", "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.md")

$de.Range("F3").Value = "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.de-de.xlf"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5834b6e726c1c7b56b89edbab07a04b41dd1997a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.de-de.xlf", "", "", "b06fe6dd-862b-45e4-b6c2-799b23eb6c96.e8ea8ac44081f12a27b9091a29751c9c3d35bb47.de-de.xlf")

$de.Range("G3").Value = "2016-02-16 14:55:46"
$de.Range("H3").Value = "Include"
